$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 206-209 (the "discharge" family entries: discharge, discharged,
# discharges, discharging). Deleting these rows shifts all subsequent rows
# up by four, which matches the new dimension A1:C513 (was A1:C517).
$ws.Range("A206:C209").EntireRow.Delete()
